# Updated cryptos list on Wed Jul 10 05:58:36 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "58.954.34"
Set-TextCell "E2" "  +2.97%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.102.50"
Set-TextCell "E3" "  +1.20%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.02%  "

# Row 5 - BNB
Set-TextCell "D5" "521.76"
Set-TextCell "E5" "  +1.55%  "

# Row 6 - Solana
Set-TextCell "D6" "143.98"
Set-TextCell "E6" "  +1.91%  "

# Row 8 - XRP
Set-TextCell "E8" "  +1.03%  "

# Row 9 - Toncoin
Set-TextCell "D9" "7.38"
Set-TextCell "E9" "  +1.74%  "

# Row 10 - Dogecoin
Set-TextCell "E10" "  +0.91%  "

# Row 11 - Cardano
Set-TextCell "E11" "  +2.66%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell "D12" "3.637.51"
Set-TextCell "E12" "  +1.30%  "

# Row 13 - TRON
Set-TextCell "E13" "  +1.29%  "

# Row 14 - Avalanche
Set-TextCell "D14" "27.08"
Set-TextCell "E14" "  +6.11%  "

# Row 15 - ShibaInu
Set-TextCell "E15" "  +0.95%  "

# Row 16 - WrappedBTC
Set-TextCell "D16" "58.959.70"
Set-TextCell "E16" "  +2.84%  "

# Row 17 & 18 swap: Polkadot <-> WrappedEther (with updated values)
Set-TextCell "B17" "Polkadot"
Set-TextCell "C17" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D17" "6.22"
Set-TextCell "E17" "  +3.16%  "

Set-TextCell "B18" "WrappedEther"
Set-TextCell "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D18" "3.107.11"
Set-TextCell "E18" "  +1.62%  "

# Row 19 - Chainlink
Set-TextCell "D19" "13.04"
Set-TextCell "E19" "  +0.36%  "

# Row 20 - Uniswap
Set-TextCell "E20" "  +1.20%  "

# Row 21 - BitcoinCash
Set-TextCell "D21" "342.84"
Set-TextCell "E21" "  +1.43%  "

# Row 23 - Polygon
Set-TextCell "E23" "  +2.03%  "

# Row 24 - Litecoin
Set-TextCell "D24" "65.83"
Set-TextCell "E24" "  +0.62%  "

# Row 25 - Kaspa
Set-TextCell "D25" "0.171"
Set-TextCell "E25" "  -0.36%  "

# Row 26 - Binance-PegBSC-USD
Set-TextCell "E26" "  -0.16%  "

# Row 27 - PEPE
Set-TextCell "D27" "0.0₃0933"
Set-TextCell "E27" "  -0.29%  "

# Row 28 - RenderToken
Set-TextCell "E28" "  +4.74%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextCell "D29" "7.27"
Set-TextCell "E29" "  +2.18%  "

# Row 30 - PancakeSwap
Set-TextCell "E30" "  +2.01%  "

# Row 31 - Fetch.AI
Set-TextCell "D31" "1.22"
Set-TextCell "E31" "  +3.82%  "

# Row 32 - EthereumClassic
Set-TextCell "D32" "21.04"
Set-TextCell "E32" "  +1.24%  "

# Row 33 - Monero
Set-TextCell "D33" "155.27"
Set-TextCell "E33" "  +0.47%  "

# Row 34 - NEARProtocol
Set-TextCell "D34" "4.67"
Set-TextCell "E34" "  +3.17%  "

# Row 35 - Aptos
Set-TextCell "D35" "6.19"
Set-TextCell "E35" "  +5.32%  "

# Row 36 - EnergySwap
Set-TextCell "D36" "27.00"
Set-TextCell "E36" "  +4.05%  "

# Row 37 - ImmutableX
Set-TextCell "E37" "  +5.60%  "

# Row 38 - Hedera
Set-TextCell "E38" "  +1.60%  "

# Row 39 - Filecoin
Set-TextCell "D39" "3.95"
Set-TextCell "E39" "  +2.69%  "

# Row 40 - RenzoRestakedETH
Set-TextCell "D40" "3.146.41"
Set-TextCell "E40" "  +1.43%  "

# Row 41 - OKB
Set-TextCell "D41" "36.89"
Set-TextCell "E41" "  -0.28%  "

# Row 42 - FirstDigitalUSD
Set-TextCell "E42" "  +0.02%  "

# Row 43 - Mantle
Set-TextCell "D43" "0.666"
Set-TextCell "E43" "  -0.47%  "

# Row 44 - Stacks
Set-TextCell "D44" "1.45"
Set-TextCell "E44" "  +5.24%  "

# Row 45 - Maker
Set-TextCell "D45" "2.284.75"
Set-TextCell "E45" "  +1.14%  "

# Row 46 - VeChain
Set-TextCell "D46" "0.0258"
Set-TextCell "E46" "  +2.66%  "

# Row 47 - InjectiveProtocol
Set-TextCell "D47" "20.96"
Set-TextCell "E47" "  +3.99%  "

# Row 48 - ONDO
Set-TextCell "D48" "0.961"
Set-TextCell "E48" "  +1.19%  "

# Row 49 - Cosmos
Set-TextCell "D49" "6.02"
Set-TextCell "E49" "  +2.78%  "

# Row 50 - SuiNetwork
Set-TextCell "D50" "0.761"
Set-TextCell "E50" "  +11.28%  "

# Row 51 - Bittensor
Set-TextCell "D51" "263.45"
Set-TextCell "E51" "  +11.46%  "
